# Apply "edit import samples & add pagination" changes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement data set: StudentID, GPA, Rank, Score
$data = @(
    @(2013011001, 3.99, 1, 50),
    @(2013011002, 3.98, 2, 51),
    @(2013011003, 3.97, 3, 52),
    @(2013011004, 3.96, 4, 53),
    @(2013011005, 3.95, 5, 54),
    @(2013011006, 3.94, 6, 55),
    @(2013011007, 3.93, 7, 56),
    @(2013011008, 3.92, 8, 57),
    @(2013011009, 3.91, 9, 58),
    @(2013011010, 3.90, 10, 59),
    @(2013011011, 3.89, 11, 60),
    @(2013011012, 3.88, 12, 61),
    @(2013011013, 3.87, 13, 62),
    @(2013011014, 3.86, 14, 63),
    @(2013011015, 3.85, 15, 64),
    @(2013011016, 3.84, 16, 65),
    @(2013011017, 3.83, 17, 66),
    @(2013011018, 3.82, 18, 67),
    @(2013011019, 3.81, 19, 68),
    @(2013011020, 3.80, 20, 69),
    @(2013011021, 3.79, 21, 70),
    @(2013011022, 3.78, 22, 71),
    @(2013011023, 3.77, 23, 72),
    @(2013011024, 3.76, 24, 73),
    @(2013011025, 3.7499999999999898, 25, 74),
    @(2013011026, 3.73999999999999, 26, 75),
    @(2013011027, 3.7299999999999902, 27, 76),
    @(2013011028, 3.71999999999999, 28, 77),
    @(2013011029, 3.7099999999999902, 29, 78)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 1
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

# Update the view: scrolled position and active selection (pagination)
$win = $excel.ActiveWindow
$win.ScrollRow = 13
$win.ScrollColumn = 1
$ws.Range("E26").Select()
